$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 360, shifting existing rows 360-445 down to 361-446
$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new data record
$ws.Range("A360").Value = 5
$ws.Range("B360").Value = "Macroferia Regional de Talca"
$ws.Range("C360").Value = "Maule"
$ws.Range("D360").Value = 44754
$ws.Range("E360").Value = 7
$ws.Range("F360").Value = 100112043
$ws.Range("G360").Value = "Pepino ensalada"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 400
$ws.Range("K360").Value = 18000
$ws.Range("L360").Value = 18000
$ws.Range("M360").Value = 18000
$ws.Range("N360").Value = "$/caja 60 unidades"
$ws.Range("O360").Value = "Región de Arica y Parinacota"
$ws.Range("P360").Value = 300
$ws.Range("Q360").Value = 60
$ws.Range("R360").Value = "Hortaliza"
